# Updates the cryptos list worksheet with refreshed price/volume data,
# matching the upstream GitHub Actions scraper commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.008.47'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '2.297.97'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.53'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.25'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.83%  '
$ws.Range("E7").Value = '  -0.67%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +1.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.24'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +7.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0791'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.96'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +4.56%  '
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("D15").Value = '2.655.58'
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("D16").Value = '2.266.36'
$ws.Range("E16").Value = '  -1.23%  '
$ws.Range("E17").Value = '  -2.16%  '
$ws.Range("D18").Value = '42.918.11'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.65'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +8.40%  '
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("E21").Value = '  +0.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.89'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.62'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.62%  '
$ws.Range("E24").Value = '  +6.33%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("E26").Value = '  -0.83%  '
$ws.Range("E27").Value = '  +1.75%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.32'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.59%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.39'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.05'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.00%  '
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("E33").Value = '  +1.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.66'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.87%  '
$ws.Range("E35").Value = '  +4.07%  '
$ws.Range("E36").Value = '  +1.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0689'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.60%  '
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.101'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.81'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.93%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  +3.13%  '
$ws.Range("E43").Value = '  -5.16%  '
$ws.Range("D44").Value = '1.979.55'
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("E45").Value = '  +2.42%  '
$ws.Range("E46").Value = '  +1.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.45'
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.46'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.69%  '
$ws.Range("E49").Value = '  +3.89%  '
$ws.Range("D50").Value = '2.521.14'
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.67'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.57%  '

